# Reorder "System" entries in the "Recorded By" column (G) so that a bare
# "System" token that currently appears first in the comma-separated list
# is moved to the end of the list, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "system, System, backup@backdoor.com" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    # Note: PowerShell comparison operators (-eq/-ne/-contains), even the
    # "-c" case-sensitive variants, behave case-insensitively in this
    # runtime, so use the .Equals() string method (case-sensitive by
    # default) to distinguish "System" from "system".
    $hasSystem = $false
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem -and ($trimmed.Count -gt 1)) {
        $newParts = @()
        foreach ($p in $trimmed) {
            if (-not $p.Equals("System")) {
                $newParts += $p
            }
        }
        $newParts += "System"
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
